# 05/08/2018 MAMATHA CHICK IN
#
# 1) "SAT Jul 21 11:20:59 IST 2018" was previously split across two runs
#    (a quirk from earlier edits); collapse it back into a single run.
# 2) Append a brand-new purchase-record block (a "chick in" entry for
#    SAT Aug 04 11:41:22 IST 2018) right after the last existing record
#    (the one ending "...Amount Received mode  - CASH AND CLEARD"),
#    ahead of the trailing blank paragraphs that close out the document.

$d = $word.ActiveDocument

# --- Part 1: merge "SAT Jul 21" + " 11:20:59 IST 2018" into one run ----------

$d.Content.Find.Execute(
    "SAT Jul 21 11:20:59 IST 2018", $true, $false, $false, $false, $false,
    $true, 1, $false, "SAT Jul 21 11:20:59 IST 2018", 2) | Out-Null

# --- Part 2: insert the new "SAT Aug 04" record block ------------------------

# Find the paragraph of the final existing record, i.e. the last
# "Amount Received mode ... - CASH AND CLEARD" paragraph in the document
# (a duplicate phrase exists earlier, so we keep walking to the last hit).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*CASH AND CLEARD*") {
        $target = $p
    }
}

$insertPoint = $d.Range($target.Range.End, $target.Range.End)

$wNs    = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rFonts = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>'

function PlainPara([string]$innerRuns) {
    return '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $rFonts + '</w:rPr></w:pPr>' + $innerRuns + '</w:p>'
}

function Run([string]$inner) {
    return '<w:r><w:rPr>' + $rFonts + '</w:rPr>' + $inner + '</w:r>'
}

function LabelValuePara([string]$label, [int]$tabCount, [string]$value) {
    $runs = Run("<w:t>$label</w:t>")
    for ($k = 1; $k -lt $tabCount; $k++) {
        $runs += Run('<w:tab/>')
    }
    $runs += Run("<w:tab/><w:t>$value</w:t>")
    return PlainPara($runs)
}

# Empty paragraph separating the previous record from this one.
$para1 = PlainPara('')

# Timestamp paragraph - two runs, exactly like the other timestamp paragraphs
# in this document (date-ish text, then the time/year text with a leading space).
$para2 = PlainPara(
    (Run('<w:t>SAT Aug 04</w:t>')) +
    (Run('<w:t xml:space="preserve"> 11:41:22 IST 2018</w:t>'))
)

$para3  = LabelValuePara 'Person Name'         4 '- KR'
$para4  = LabelValuePara 'Bill number'         4 '- 7061'
$para5  = PlainPara((Run('<w:t>---------------------------------------------------------------</w:t>')))
$para6  = LabelValuePara 'Item Name'           4 '- POTATO'
$para7  = LabelValuePara 'Number of Pockets'   3 '- 3'
$para8  = LabelValuePara 'Number of KGs'       3 '- 162'
$para9  = LabelValuePara 'Rate'                5 '- 20'
$para10 = LabelValuePara 'Total Price'         4 '- 3240.0'

# "Amount balance" is bold.
$boldRuns = '<w:r><w:rPr>' + $rFonts + '<w:b/></w:rPr><w:t>Amount balance</w:t></w:r>' +
            '<w:r><w:rPr>' + $rFonts + '<w:b/></w:rPr><w:tab/></w:r>' +
            '<w:r><w:rPr>' + $rFonts + '<w:b/></w:rPr><w:tab/></w:r>' +
            '<w:r><w:rPr>' + $rFonts + '<w:b/></w:rPr><w:tab/><w:t>- 3240.0</w:t></w:r>'
$para11 = '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $rFonts + '<w:b/></w:rPr></w:pPr>' + $boldRuns + '</w:p>'

# Two trailing blank paragraphs, as in the diff.
$para12 = PlainPara('')
$para13 = PlainPara('')

$newBlockXml = '<w:p ' + $wNs + '>' +
    ($para1 -replace '^<w:p>', '') +
    $para2 + $para3 + $para4 + $para5 + $para6 + $para7 + $para8 + $para9 +
    $para10 + $para11 + $para12 + $para13

$insertPoint.InsertXML($newBlockXml)
